# The "Functions" sheet has a set of cells in columns C/D that hold a
# placeholder test value. Almost all of them already read "Y-Test", but
# four cells (C18, D18, C19, C20) still held the old placeholder "X".
# Bring them in line with the rest of the column.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Functions")

$ws.Range("C18").Value = "Y-Test"
$ws.Range("D18").Value = "Y-Test"
$ws.Range("C19").Value = "Y-Test"
$ws.Range("C20").Value = "Y-Test"
